$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.004
$ws.Range("E2").Value = 0.834
$ws.Range("G2").Value = 0
$ws.Range("I2").Value = 0.004

$ws.Range("D3").Value = 0.007
$ws.Range("E3").Value = 0.795
$ws.Range("F3").Value = 0.001
$ws.Range("H3").Value = 0.001
$ws.Range("I3").Value = 0.006

$ws.Range("D4").Value = 0.024
$ws.Range("E4").Value = 0.764
$ws.Range("F4").Value = 0.003
$ws.Range("G4").Value = 0.012
$ws.Range("H4").Value = 0.008
$ws.Range("I4").Value = 0.022

$ws.Range("D5").Value = 0.064
$ws.Range("E5").Value = 0.702
$ws.Range("F5").Value = 0.022
$ws.Range("G5").Value = 0.042
$ws.Range("H5").Value = 0.036
$ws.Range("I5").Value = 0.057

$ws.Range("D6").Value = 0.22
$ws.Range("E6").Value = 0.6
$ws.Range("F6").Value = 0.155
$ws.Range("G6").Value = 0.179
$ws.Range("I6").Value = 0.214

$ws.Range("D7").Value = 0.523
$ws.Range("E7").Value = 0.493
$ws.Range("F7").Value = 0.472
$ws.Range("G7").Value = 0.506
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = 0.538

$ws.Range("D8").Value = 0.852
$ws.Range("E8").Value = 0.37
$ws.Range("F8").Value = 0.902
$ws.Range("G8").Value = 0.882
$ws.Range("H8").Value = 0.892
$ws.Range("I8").Value = 0.868

$ws.Range("D9").Value = 0.995
$ws.Range("E9").Value = 0.214
$ws.Range("F9").Value = 0.998
$ws.Range("G9").Value = 0.998
$ws.Range("H9").Value = 0.998
$ws.Range("I9").Value = 0.998

$ws.Range("E10").Value = 0.078

$ws.Range("E11").Value = 0.007

$ws.Range("E13").Value = 0.493
$ws.Range("F13").Value = 0.841
$ws.Range("H13").Value = 0.021

$ws.Range("E14").Value = 0.531
$ws.Range("F14").Value = 0.818
$ws.Range("H14").Value = 0.047

$ws.Range("F15").Value = 0.752
$ws.Range("H15").Value = 0.064

$ws.Range("E16").Value = 0.522
$ws.Range("F16").Value = 0.694
$ws.Range("H16").Value = 0.162

$ws.Range("D17").Value = 0.997
$ws.Range("E17").Value = 0.501
$ws.Range("F17").Value = 0.62
$ws.Range("G17").Value = 0.98
$ws.Range("H17").Value = 0.273
$ws.Range("I17").Value = 0.994

$ws.Range("D18").Value = 0.513
$ws.Range("E18").Value = 0.511
$ws.Range("F18").Value = 0.483
$ws.Range("G18").Value = 0.493
$ws.Range("H18").Value = 0.484
$ws.Range("I18").Value = 0.521

$ws.Range("D19").Value = 0.003
$ws.Range("E19").Value = 0.49
$ws.Range("F19").Value = 0.383
$ws.Range("G19").Value = 0.022
$ws.Range("H19").Value = 0.749
$ws.Range("I19").Value = 0.006

$ws.Range("E20").Value = 0.475
$ws.Range("F20").Value = 0.27
$ws.Range("H20").Value = 0.939

$ws.Range("E21").Value = 0.455
$ws.Range("F21").Value = 0.158
$ws.Range("H21").Value = 0.994

$ws.Range("E22").Value = 0.452
$ws.Range("F22").Value = 0.069

$ws.Range("E23").Value = 0.416
$ws.Range("F23").Value = 0.025

$ws.Range("D24").Value = 0.006
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.805
$ws.Range("G24").Value = 0.39
$ws.Range("H24").Value = 0.614
$ws.Range("I24").Value = 0.016

$ws.Range("D25").Value = 0.01
$ws.Range("E25").Value = 0.992
$ws.Range("F25").Value = 0.772
$ws.Range("G25").Value = 0.395
$ws.Range("H25").Value = 0.588
$ws.Range("I25").Value = 0.022

$ws.Range("D26").Value = 0.028
$ws.Range("E26").Value = 0.978
$ws.Range("F26").Value = 0.728
$ws.Range("H26").Value = 0.591
$ws.Range("I26").Value = 0.055

$ws.Range("D27").Value = 0.085
$ws.Range("E27").Value = 0.928
$ws.Range("F27").Value = 0.681
$ws.Range("G27").Value = 0.433
$ws.Range("H27").Value = 0.583
$ws.Range("I27").Value = 0.109

$ws.Range("D28").Value = 0.225
$ws.Range("E28").Value = 0.758
$ws.Range("F28").Value = 0.576
$ws.Range("G28").Value = 0.464
$ws.Range("H28").Value = 0.567
$ws.Range("I28").Value = 0.255

$ws.Range("D29").Value = 0.528
$ws.Range("E29").Value = 0.486
$ws.Range("F29").Value = 0.496
$ws.Range("G29").Value = 0.497
$ws.Range("H29").Value = 0.496
$ws.Range("I29").Value = 0.529

$ws.Range("D30").Value = 0.839
$ws.Range("E30").Value = 0.18
$ws.Range("F30").Value = 0.382
$ws.Range("G30").Value = 0.553
$ws.Range("H30").Value = 0.445
$ws.Range("I30").Value = 0.793

$ws.Range("D31").Value = 0.983
$ws.Range("E31").Value = 0.02
$ws.Range("F31").Value = 0.29
$ws.Range("G31").Value = 0.584
$ws.Range("H31").Value = 0.412
$ws.Range("I31").Value = 0.97

$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0.162
$ws.Range("G32").Value = 0.646
$ws.Range("H32").Value = 0.35
$ws.Range("I32").Value = 0.999

$ws.Range("F33").Value = 0.099
$ws.Range("G33").Value = 0.7
$ws.Range("H33").Value = 0.289

$ws.Range("F34").Value = 0.032
$ws.Range("G34").Value = 0.795
$ws.Range("H34").Value = 0.187

$ws.Range("E35").Value = 0.001
$ws.Range("F35").Value = 0.975
$ws.Range("G35").Value = 0.006
$ws.Range("I35").Value = 0.016

$ws.Range("D36").Value = 0.021
$ws.Range("E36").Value = 0.001
$ws.Range("F36").Value = 0.959
$ws.Range("G36").Value = 0.01
$ws.Range("I36").Value = 0.027

$ws.Range("D37").Value = 0.075
$ws.Range("E37").Value = 0.01
$ws.Range("F37").Value = 0.924
$ws.Range("G37").Value = 0.029
$ws.Range("I37").Value = 0.063

$ws.Range("D38").Value = 0.112
$ws.Range("E38").Value = 0.03
$ws.Range("F38").Value = 0.854
$ws.Range("G38").Value = 0.066
$ws.Range("H38").Value = 0.999
$ws.Range("I38").Value = 0.11

$ws.Range("D39").Value = 0.255
$ws.Range("E39").Value = 0.136
$ws.Range("F39").Value = 0.709
$ws.Range("G39").Value = 0.191
$ws.Range("H39").Value = 0.935
$ws.Range("I39").Value = 0.26

$ws.Range("D40").Value = 0.496
$ws.Range("E40").Value = 0.491
$ws.Range("F40").Value = 0.481
$ws.Range("G40").Value = 0.497
$ws.Range("H40").Value = 0.517
$ws.Range("I40").Value = 0.506

$ws.Range("D41").Value = 0.819
$ws.Range("E41").Value = 0.906
$ws.Range("F41").Value = 0.213
$ws.Range("G41").Value = 0.866
$ws.Range("H41").Value = 0.025
$ws.Range("I41").Value = 0.831

$ws.Range("D42").Value = 0.981
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 0.033
$ws.Range("G42").Value = 0.994
$ws.Range("I42").Value = 0.984

$ws.Range("F43").Value = 0.004
